# Re-run SGNN to annotate dialog acts following clean up work to the
# original transcripts. Updates DAMSLTag (col I) / DialogAct (col J)
# values for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    22  = @("sd", "Statement-non-opinion")
    53  = @("aa", "Agree/Accept")
    58  = @("aa", "Agree/Accept")
    61  = @("sd", "Statement-non-opinion")
    65  = @("sd", "Statement-non-opinion")
    69  = @("sd", "Statement-non-opinion")
    73  = @("sd", "Statement-non-opinion")
    87  = @("sd", "Statement-non-opinion")
    93  = @("qy", "Yes-No-Question")
    103 = @("sd", "Statement-non-opinion")
    104 = @("aa", "Agree/Accept")
    108 = @("aa", "Agree/Accept")
    111 = @("sd", "Statement-non-opinion")
    112 = @("sv", "Statement-opinion")
    116 = @("sv", "Statement-opinion")
    123 = @("sd", "Statement-non-opinion")
    135 = @("sv", "Statement-opinion")
    142 = @("sd", "Statement-non-opinion")
    158 = @("sv", "Statement-opinion")
    159 = @("aa", "Agree/Accept")
    163 = @("aa", "Agree/Accept")
    176 = @("aa", "Agree/Accept")
    197 = @("sd", "Statement-non-opinion")
    216 = @("aa", "Agree/Accept")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
